$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Driver Summary")

$ws.Range("C3").Value = 377
$ws.Range("C4").Value = 120
$ws.Range("D4").Value = 96.7
$ws.Range("C5").Value = 2223
$ws.Range("C6").Value = 2720
